$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the order of the "Recorded By" entries from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# across the whole used range of the "Recorded By" column (G).
$rng = $ws.UsedRange
$rng.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com", 1, 1, $false, $false, $false)
